$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update direct values
$ws.Range("C2").Value = 2246
$ws.Range("D2").Value = 1261

# Update formula in G2
$ws.Range("G2").Formula = "=98"

# Update the active cell selection on the sheet to F2
$ws.Range("F2").Select()
